# Extend the GT4500 test-case sheet with a 5th test case and tweak a
# couple of existing descriptions, per the "Extended tests for GT4500 due
# to lack of minimum coverage" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Tweak three existing descriptions (2nd and 4th test cases).
# ---------------------------------------------------------------------
$ws.Range("B9").Value = "Az egyik tároló tartalmaz 2 torpedót, a másik üres (mindkét tárolóval tesztelve)"
$ws.Range("B11").Value = "A második és a harmadik, valamint az ötödik és hatodik parancs eredménye SUCCESS"
$ws.Range("B18").Value = "A GT4500 tud egyszerre tüzelni mindkét tárolóból, ha az egyik üres, akkor a másikból tüzel, ha mindkettő üres, akkor nem tüzel"

# ---------------------------------------------------------------------
# 2) Fill in the brand-new 5th test case in rows 22-26 (previously blank
#    filler rows), re-using the same layout as the other test cases.
# ---------------------------------------------------------------------
$ws.Range("A22").Value = "5. teszteset"
$ws.Range("A23").Value = "Tesztelt követelmény:"
$ws.Range("B23").Value = "A GT4500 csak a SINGLE és ALL tüzelési módokat ismeri"
$ws.Range("A24").Value = "Előfeltételek:"
$ws.Range("B24").Value = "Van egy GT4500 hajó"
$ws.Range("A25").Value = "A teszt lépései"
$ws.Range("B25").Value = "Hibás tüzelési módot adunk meg a GT4500-nak (pl. RANDOM)"
$ws.Range("A26").Value = "Elvárt eredmény:"
$ws.Range("B26").Value = "A második parancs eredménye Unknown firing mode: 'RANDOM'"

# Give the new section header row (22) the same look as the other
# section header rows (2, 7, 12, 17): centered, top aligned, wrapped.
$headerRow = $ws.Range("A22:B22")
$headerRow.HorizontalAlignment = -4108   # xlCenter
$headerRow.VerticalAlignment = -4160     # xlTop
$headerRow.WrapText = $true
$headerRow.Merge()

# The four detail rows of the new test case (23-26) match the plain
# "wrap + top, general alignment" look already used by every other
# detail row.
$detailRows = $ws.Range("A23:B26")
$detailRows.HorizontalAlignment = 1      # xlGeneral
$detailRows.VerticalAlignment = -4160    # xlTop
$detailRows.WrapText = $true

# ---------------------------------------------------------------------
# 3) Cosmetic view updates: wider columns and the new selection.
# ---------------------------------------------------------------------
$ws.Range("A1:B1").EntireColumn.ColumnWidth = 119.83

$ws.Range("B24").Select()
